$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defaults_20130301")

# Row 5 (SAL) - recalc using TDS: H5 was a formula (G5+I5)/2, now a plain stepped value.
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 75
$ws.Range("L5").Value = 80
$ws.Range("M5").Value = 80

# Row 22 - L/M updated
$ws.Range("L22").Value = 333
$ws.Range("M22").Value = 333

# New column P on row 24
$ws.Range("P24").Formula = "=80 * (50/12)"

# Move selection to reflect the author's last active cell
$ws.Range("G8").Select() | Out-Null
